$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3192.6902
$ws.Range("I15").Value = 3192.6902
$ws.Range("K15").Value = 9578.070599999999
$ws.Range("M15").Value = -9409.070599999999
$ws.Range("H33").Value = 311
$ws.Range("J33").Value = 124.5
$ws.Range("L33").Value = 124.5
$ws.Range("N33").Value = -582.5
$ws.Range("H112").Value = 2389
$ws.Range("J112").Value = 2389
$ws.Range("L112").Value = 7167
$ws.Range("N112").Value = -9383
$ws.Range("H125").Value = 1120
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1120
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10080
$ws.Range("N125").Value = -15000
$ws.Range("H137").Value = 10021126
$ws.Range("I137").Value = 16668141
$ws.Range("K137").Value = 50004423
$ws.Range("M137").Value = -50001873
$ws.Range("M125").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4339.2905
$ws.Range("I61").Value = 2937.7273
$ws.Range("K61").Value = 2937.7273
$ws.Range("M61").Value = -2725.7273
$ws.Range("H63").Value = 5111.6665
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 5111.6665
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H132").Value = 2642.348
$ws.Range("I132").Value = 1879.6875
$ws.Range("J132").Value = 4385.5713
$ws.Range("K132").Value = 5639.0625
$ws.Range("L132").Value = 13156.7139
$ws.Range("M132").Value = -3109.0625
$ws.Range("N132").Value = -18216.7139
$ws.Range("H136").Value = 4339.2905
$ws.Range("I136").Value = 2937.7273
$ws.Range("K136").Value = 8813.1819
$ws.Range("M136").Value = -6263.1819
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17860706
$ws.Range("I20").Value = 26320110
$ws.Range("J20").Value = 1961.5555
$ws.Range("K20").Value = 26320110
$ws.Range("L20").Value = 1961.5555
$ws.Range("M20").Value = -26319863
$ws.Range("N20").Value = -2455.5555
$ws.Range("H82").Value = 51054.3
$ws.Range("I82").Value = 24999.334
$ws.Range("J82").Value = 90136.75
$ws.Range("K82").Value = 24999.334
$ws.Range("L82").Value = 90136.75
$ws.Range("M82").Value = -24616.334
$ws.Range("N82").Value = -90902.75
$ws.Range("H85").Value = 51054.3
$ws.Range("I85").Value = 24999.334
$ws.Range("J85").Value = 90136.75
$ws.Range("K85").Value = 24999.334
$ws.Range("L85").Value = 90136.75
$ws.Range("M85").Value = -23673.334
$ws.Range("N85").Value = -92788.75
$ws.Range("H134").Value = 4356.9365
$ws.Range("I134").Value = 4238.314
$ws.Range("J134").Value = 4861.0835
$ws.Range("K134").Value = 12714.942
$ws.Range("L134").Value = 14583.2505
$ws.Range("M134").Value = -10179.942
$ws.Range("N134").Value = -19653.2505

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5275.5864
$ws.Range("I31").Value = 4109.5
$ws.Range("K31").Value = 4109.5
$ws.Range("M31").Value = -3814.5
$ws.Range("H34").Value = 5275.5864
$ws.Range("I34").Value = 4109.5
$ws.Range("K34").Value = 4109.5
$ws.Range("M34").Value = -3907.5
$ws.Range("H68").Value = 88473.5
$ws.Range("J68").Value = 88473.5
$ws.Range("L68").Value = 88473.5
$ws.Range("N68").Value = -89971.5
$ws.Range("H71").Value = 88473.5
$ws.Range("J71").Value = 88473.5
$ws.Range("L71").Value = 265420.5
$ws.Range("N71").Value = -272908.5
$ws.Range("H132").Value = 2734.5789
$ws.Range("I132").Value = 2121.6
$ws.Range("J132").Value = 3913.3845
$ws.Range("K132").Value = 6364.799999999999
$ws.Range("L132").Value = 11740.1535
$ws.Range("M132").Value = -3834.799999999999
$ws.Range("N132").Value = -16800.1535

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5931.048
$ws.Range("I56").Value = 5931.048
$ws.Range("K56").Value = 5931.048
$ws.Range("M56").Value = -5401.048
$ws.Range("H98").Value = 2251.5
$ws.Range("I98").Value = 1648
$ws.Range("J98").Value = 2855
$ws.Range("K98").Value = 4944
$ws.Range("L98").Value = 8565
$ws.Range("M98").Value = -3446
$ws.Range("N98").Value = -11561

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 70825.89999999999
$ws.Range("I70").Value = 103984.95
$ws.Range("J70").Value = 4507.8
$ws.Range("K70").Value = 103984.95
$ws.Range("L70").Value = 4507.8
$ws.Range("M70").Value = -103714.95
$ws.Range("N70").Value = -5047.8
$ws.Range("H73").Value = 70825.89999999999
$ws.Range("I73").Value = 103984.95
$ws.Range("J73").Value = 4507.8
$ws.Range("K73").Value = 103984.95
$ws.Range("L73").Value = 4507.8
$ws.Range("M73").Value = -103048.95
$ws.Range("N73").Value = -6379.8
$ws.Range("H80").Value = 142859280
$ws.Range("J80").Value = 1996.6666
$ws.Range("L80").Value = 1996.6666
$ws.Range("N80").Value = -3992.6666
$ws.Range("H83").Value = 142859280
$ws.Range("J83").Value = 1996.6666
$ws.Range("L83").Value = 9983.333000000001
$ws.Range("N83").Value = -19967.333
$ws.Range("H102").Value = 8087.6875
$ws.Range("I102").Value = 1376.6666
$ws.Range("K102").Value = 1376.6666
$ws.Range("M102").Value = 245.3334
$ws.Range("H112").Value = 92823.25
$ws.Range("J112").Value = 92823.25
$ws.Range("L112").Value = 92823.25
$ws.Range("N112").Value = -95039.25
$ws.Range("H132").Value = 3237.6
$ws.Range("I132").Value = 1800
$ws.Range("K132").Value = 5400
$ws.Range("M132").Value = -2870

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2764.3076
$ws.Range("I7").Value = 2829.6667
$ws.Range("K7").Value = 2829.6667
$ws.Range("M7").Value = -2717.6667
$ws.Range("H100").Value = 3699.4167
$ws.Range("J100").Value = 2548
$ws.Range("L100").Value = 2548
$ws.Range("N100").Value = -3630
$ws.Range("H122").Value = 4550.625
$ws.Range("I122").Value = 4858.143
$ws.Range("K122").Value = 14574.429
$ws.Range("M122").Value = -12124.429
$ws.Range("H126").Value = 2764.3076
$ws.Range("I126").Value = 2829.6667
$ws.Range("K126").Value = 8489.000100000001
$ws.Range("M126").Value = -6019.000100000001
$ws.Range("H132").Value = 4772.148
$ws.Range("I132").Value = 2781.2632
$ws.Range("J132").Value = 9500.5
$ws.Range("K132").Value = 8343.7896
$ws.Range("L132").Value = 28501.5
$ws.Range("M132").Value = -5813.7896
$ws.Range("N132").Value = -33561.5
$ws.Range("H136").Value = 6009.5557
$ws.Range("I136").Value = 7523.5
$ws.Range("K136").Value = 22570.5
$ws.Range("M136").Value = -20020.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H132").Value = 2377.3125
$ws.Range("I132").Value = 2118.5
$ws.Range("J132").Value = 3498.8333
$ws.Range("K132").Value = 6355.5
$ws.Range("L132").Value = 10496.4999
$ws.Range("M132").Value = -3825.5
$ws.Range("N132").Value = -15556.4999
$ws.Range("N105").ClearContents()
$ws.Range("N93").ClearContents()
